$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header text between K1 and L1
# (K1 was "79 required function fill in internal", L1 was "80 . preceeded by :::")
$ws.Range("K1").Value = "80`n. preceeded by :::"
$ws.Range("L1").Value = "79`nrequired function fill in internal"

# Move the "x" mark in row 2 from L2 to K2 (follows the header that moved to K1)
$ws.Range("K2").Value = "x"
$ws.Range("L2").Clear()

# Add a new "x" mark in row 9
$ws.Range("M9").Value = "x"

# Fill in the remaining "x" marks across row 18
$ws.Range("C18").Value = "x"
$ws.Range("E18").Value = "x"
$ws.Range("G18").Value = "x"
$ws.Range("H18").Value = "x"
$ws.Range("I18").Value = "x"
$ws.Range("J18").Value = "x"
$ws.Range("K18").Value = "x"
$ws.Range("L18").Value = "x"
$ws.Range("M18").Value = "x"

# Update the selected cell to match the new cursor position
$ws.Range("M12").Select()
